$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 is the "SamplesTab" row: B3 keeps the Sample ID query, now with the
# Tumor / Analyte Type columns removed from the SELECT list
$ws.Range("B3").Value = "SELECT`n    DISTINCT (smp.sample_id) AS `"Sample ID`",`n    sp.participant_id AS `"Participant ID`", `n    s.study_name AS `"Study Name`",`n    s.phs_accession AS Accession`nFROM `n    df_participant sp`nJOIN `n    df_study s ON sp.`"study.phs_accession`" = s.phs_accession`nJOIN `n    df_sample smp ON smp.`"participant.study_participant_id`" = sp.study_participant_id`nJOIN`n    df_diagnosis d ON d.`"participant.study_participant_id`" = sp.study_participant_id`nJOIN`n    df_program p ON p.program_acronym = s.`"program.program_acronym`"`nJOIN`n    df_file f1 ON f1.`"sample.sample_id`" = smp.sample_id`nJOIN`n    df_genomic_info gi ON gi.`"file.file_id`" = f1.file_id`nWHERE `n   s.phs_accession = 'phs001437' AND gi.library_layout = 'Paired-End'`nORDER BY `n    smp.sample_id ASC`nLIMIT 100;"

# Row 4 is the "FilesTab" row: B4 keeps the File Name query (unchanged text)
$ws.Range("B4").Value = "SELECT`n    f1.file_name AS `"File Name`",`n    s.study_name AS `"Study Name`",`n    s.phs_accession AS `"Accession`",`n    sp.participant_id AS `"Participant Id`",`n    COALESCE((`n        SELECT`n            REPLACE(GROUP_CONCAT(CASE WHEN rn <= 5 THEN smp.sample_id ELSE NULL END, ', '), ', , ', ', ') ||`n            CASE WHEN MAX(rn) > 5 THEN ', ...' ELSE '' END`n        FROM (`n            SELECT`n                smp.sample_id,`n                ROW_NUMBER() OVER (ORDER BY smp.sample_id) AS rn`n            FROM df_sample smp`n            WHERE smp.`"participant.study_participant_id`" = sp.study_participant_id`n        ) smp`n    ), '') AS `"Sample Id`",`n    f1.file_type AS `"File Type`",`n    gi.library_strategy AS `"Library Strategy`"`nFROM `n    df_study s`nINNER JOIN `n    df_participant sp ON sp.`"study.phs_accession`" = s.phs_accession`nINNER JOIN  `n    df_sample smp ON smp.`"participant.study_participant_id`" = sp.study_participant_id`nINNER JOIN `n    df_file f1 ON f1.`"sample.sample_id`" = smp.sample_id`nINNER JOIN`n    df_genomic_info gi ON gi.`"file.file_id`" = f1.file_id`nINNER JOIN`n    df_diagnosis d ON d.`"participant.study_participant_id`" = sp.study_participant_id`nINNER JOIN`n    df_program p ON p.program_acronym = s.`"program.program_acronym`"`nWHERE `n    s.phs_accession = 'phs001437' AND gi.library_layout = 'Paired-End'`nGROUP BY`n    f1.file_name,`n    s.study_name,`n    s.phs_accession,`n    sp.participant_id,`n    f1.file_type,`n    gi.library_strategy`nORDER BY `n    f1.file_name ASC`nLIMIT 100;"

# The TSV/Web filename columns (D/E) are no longer populated for rows 3 and 4
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# Update the active selection to match the authored state
$ws.Range("D3:E4").Select()
